$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary updates ---
# Valor Mora (total) 790666 -> 390666
$ws.Range("E11").Value = 390666

# Cant. Trabajadores 5 -> 2
$ws.Range("C13").Value = 2

# --- Worker / period detail rows ---
# Row 16 (Walter Enrique Jimenez Diaz, period 2210) stays unchanged.

# Row 17: was GEINER DE LA CRUZ RIQUET / 2309 -> now CESAR DAVID ORTEGA MORALES / 2210
$ws.Range("C17").Value = "1049826190"
$ws.Range("D17").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E17").Value = "2210"

# Row 18: GEINER / 2308 -> CESAR / 2303
$ws.Range("C18").Value = "1049826190"
$ws.Range("D18").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E18").Value = "2303"

# Row 19: GEINER / 2307 -> CESAR / 2304
$ws.Range("C19").Value = "1049826190"
$ws.Range("D19").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E19").Value = "2304"

# Row 20: GEINER / 2306 -> CESAR / 2305
$ws.Range("C20").Value = "1049826190"
$ws.Range("D20").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E20").Value = "2305"

# Row 21: GEINER / 2305 -> CESAR / 2306
$ws.Range("C21").Value = "1049826190"
$ws.Range("D21").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E21").Value = "2306"

# Row 22: GEINER / 2304 -> CESAR / 2307
$ws.Range("C22").Value = "1049826190"
$ws.Range("D22").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E22").Value = "2307"

# Row 23: GEINER / 2303 -> CESAR / 2308
$ws.Range("C23").Value = "1049826190"
$ws.Range("D23").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E23").Value = "2308"

# Row 24: GEINER / 2210 -> CESAR / 2309
$ws.Range("C24").Value = "1049826190"
$ws.Range("D24").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E24").Value = "2309"

# Row 25 already holds CESAR DAVID ORTEGA MORALES / 2310 with the trailing balance values
# (F25=30666, G25=1000000); it becomes the new last row of the table and needs the
# "closing" border/style treatment that used to belong to the old last row (35).
$ws.Range("C25").Value = "1049826190"
$ws.Range("D25").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E25").Value = "2310"

$ws.Range("B25:J25").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws.Range("B25:J25").Borders.Item(9).Weight = 2

# Copy the "closing row" style that used to live on row 35 onto row 25, then remove the
# now-obsolete rows 26-35 (GEINER/ELBIS/AMAURY legacy rows that are no longer part of the
# account statement).
$ws.Range("B35:J35").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C25").Value = "1049826190"
$ws.Range("D25").Value = "CESAR DAVID ORTEGA MORALES"
$ws.Range("E25").Value = "2310"

$ws.Rows("26:35").Delete()

$excel.CutCopyMode = 0
